$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new series values: P1=14, Q1=15
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the formatting (bold/border/center style) from O1 onto the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update data rows 2-25: swap values in columns I/K and M/O, and add new
# columns P and Q (both filled with 2) to extend the data table.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2 (new column)
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2 (new column)
}
